$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.900.43'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.624.78'
$ws.Range("E3").Value = '  +0.99%  '

$ws.Range("E4").Value = '  -0.42%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.48'
$ws.Range("E5").Value = '  +0.78%  '

$ws.Range("E6").Value = '  +0.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.40%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '30.01'
$ws.Range("E8").Value = '  +11.35%  '

$ws.Range("E9").Value = '  +2.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0612'
$ws.Range("E10").Value = '  +1.60%  '

$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.855.48'
$ws.Range("E12").Value = '  +0.95%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.618.26'
$ws.Range("E13").Value = '  +0.56%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.571'
$ws.Range("E14").Value = '  +6.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.91'
$ws.Range("E15").Value = '  +5.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.932.56'
$ws.Range("E16").Value = '  +1.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.86'
$ws.Range("E17").Value = '  +16.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.70'
$ws.Range("E18").Value = '  +1.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.46'
$ws.Range("E19").Value = '  +1.53%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0706'
$ws.Range("E20").Value = '  +1.61%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.997'
$ws.Range("E21").Value = '  -0.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("E22").Value = '  +3.26%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.64'
$ws.Range("E23").Value = '  +4.23%  '

$ws.Range("E24").Value = '  +1.95%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.30'
$ws.Range("E25").Value = '  +1.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.71'
$ws.Range("E26").Value = '  +2.52%  '

$ws.Range("E27").Value = '  +2.21%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.62'
$ws.Range("E28").Value = '  +2.99%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0489'
$ws.Range("E30").Value = '  +2.95%  '

$ws.Range("E31").Value = '  +5.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  +3.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.23'
$ws.Range("E33").Value = '  +3.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.426.61'
$ws.Range("E34").Value = '  +0.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.65'
$ws.Range("E35").Value = '  +6.88%  '

$ws.Range("E36").Value = '  -0.09%  '

$ws.Range("E37").Value = '  +1.86%  '

$ws.Range("E38").Value = '  -0.70%  '

$ws.Range("E39").Value = '  +2.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.560'
$ws.Range("E40").Value = '  +3.48%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.00'
$ws.Range("E41").Value = '  +0.44%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.834'
$ws.Range("E42").Value = '  +4.24%  '

$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0501'
$ws.Range("E43").Value = '  +1.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '54.38'
$ws.Range("E44").Value = '  +0.46%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.35'
$ws.Range("E45").Value = '  +5.06%  '

$ws.Range("E46").Value = '  +16.42%  '

$ws.Range("E47").Value = '  -0.46%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.42'
$ws.Range("E48").Value = '  +2.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.764.96'
$ws.Range("E49").Value = '  +0.97%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '88.73'
$ws.Range("E50").Value = '  +2.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0108'
$ws.Range("E51").Value = '  +3.62%  '
